$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the date-column style (s="2") from A660 down through A672
$ws.Range("A660").Copy()
$ws.Range("A661:A672").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$data = @(
    @(661, 45129, 4550.16015625, 4555, 4535.7900390625, 4536.33984375, 4536.33984375, 3570190000),
    @(662, 45130, 4550.16015625, 4555, 4535.7900390625, 4536.33984375, 4536.33984375, 3570190000),
    @(663, 45131, 4543.39013671875, 4563.41015625, 4541.2900390625, 4554.64013671875, 4554.64013671875, 3856250000),
    @(664, 45132, 4555.18994140625, 4580.6201171875, 4552.419921875, 4567.4599609375, 4567.4599609375, 3812470000),
    @(665, 45133, 4558.9599609375, 4582.47021484375, 4547.580078125, 4566.75, 4566.75, 3990290000),
    @(666, 45134, 4598.259765625, 4607.06982421875, 4528.56005859375, 4537.41015625, 4537.41015625, 4553210000),
    @(667, 45135, 4565.75, 4590.16015625, 4564.009765625, 4582.22998046875, 4582.22998046875, 3981010000),
    @(668, 45136, 4565.75, 4590.16015625, 4564.009765625, 4582.22998046875, 4582.22998046875, 3981010000),
    @(669, 45137, 4565.75, 4590.16015625, 4564.009765625, 4582.22998046875, 4582.22998046875, 3981010000),
    @(670, 45138, 4584.81982421875, 4594.22021484375, 4573.14013671875, 4588.9599609375, 4588.9599609375, 4503600000),
    @(671, 45139, 4578.830078125, 4584.6201171875, 4567.52978515625, 4576.72998046875, 4576.72998046875, 4042370000),
    @(672, 45140, 4550.93017578125, 4550.93017578125, 4505.75, 4513.39013671875, 4513.39013671875, 4270710000)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value2 = $row[1]
    $ws.Cells.Item($r, 2).Value2 = $row[2]
    $ws.Cells.Item($r, 3).Value2 = $row[3]
    $ws.Cells.Item($r, 4).Value2 = $row[4]
    $ws.Cells.Item($r, 5).Value2 = $row[5]
    $ws.Cells.Item($r, 6).Value2 = $row[6]
    $ws.Cells.Item($r, 7).Value2 = $row[7]
}
